$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (style, number format, borders, etc.) of the last
# existing data row (328) down onto the new rows (329-343) so the new
# cells pick up the same cellXf (date style on column A, etc.) instead of
# creating new style entries.
$srcRow = $ws.Range("A328:D328")
$destRows = $ws.Range("A329:D343")
$srcRow.Copy($destRows)

# New daily data appended through 2021-08-09 (serials 44403-44417).
$data = @(
    @(329, 44403, 0, 1, 26.76659528907923),
    @(330, 44404, 0, 1, 26.76659528907923),
    @(331, 44405, 0, 1, 26.76659528907923),
    @(332, 44406, 0, 1, 26.76659528907923),
    @(333, 44407, 0, 1, 26.76659528907923),
    @(334, 44408, 0, 1, 26.76659528907923),
    @(335, 44409, 0, 0, 0),
    @(336, 44410, 0, 0, 0),
    @(337, 44411, 0, 0, 0),
    @(338, 44412, 0, 0, 0),
    @(339, 44413, 1, 1, 26.76659528907923),
    @(340, 44414, 1, 2, 53.53319057815846),
    @(341, 44415, 2, 4, 107.0663811563169),
    @(342, 44416, 1, 5, 133.8329764453961),
    @(343, 44417, 0, 5, 133.8329764453961)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}
